# Add new test case (TC-004-style row) to the VenusSystem test sheet.
# Row 5 mirrors row 4 (same test data), except the ExpectedResult
# (column N) is "ERROR" instead of "VALIDATION".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 4's values + formatting into the new row 5.
$ws.Range("A4:N4").Copy($ws.Range("A5:N5"))

# The new test case expects an ERROR result instead of VALIDATION.
$ws.Range("N5").Value = "ERROR"

# Match the author's final selection/active cell.
$ws.Range("N5").Select()
